$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 104
$ws.Cells.Item(104, 1).Value = 9
$ws.Cells.Item(104, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(104, 3).Value = 'Metropolitana'
$ws.Cells.Item(104, 4).Value = 44776
$ws.Cells.Item(104, 5).Value = 13
$ws.Cells.Item(104, 6).Value = 100112022
$ws.Cells.Item(104, 7).Value = 'Arveja Verde'
$ws.Cells.Item(104, 8).Value = 'Perfection'
$ws.Cells.Item(104, 9).Value = 'Primera'
$ws.Cells.Item(104, 10).Value = 25
$ws.Cells.Item(104, 11).Value = 42000
$ws.Cells.Item(104, 12).Value = 42000
$ws.Cells.Item(104, 13).Value = 42000
$ws.Cells.Item(104, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(104, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(104, 16).Value = 1680
$ws.Cells.Item(104, 17).Value = 25
$ws.Cells.Item(104, 18).Value = 'Hortaliza'

# Row 105
$ws.Cells.Item(105, 1).Value = 9
$ws.Cells.Item(105, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(105, 3).Value = 'Metropolitana'
$ws.Cells.Item(105, 4).Value = 44301
$ws.Cells.Item(105, 5).Value = 13
$ws.Cells.Item(105, 6).Value = 100112022
$ws.Cells.Item(105, 7).Value = 'Arveja Verde'
$ws.Cells.Item(105, 8).Value = 'Perfection'
$ws.Cells.Item(105, 9).Value = 'Primera'
$ws.Cells.Item(105, 10).Value = 30
$ws.Cells.Item(105, 11).Value = 32000
$ws.Cells.Item(105, 12).Value = 32000
$ws.Cells.Item(105, 13).Value = 32000
$ws.Cells.Item(105, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(105, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(105, 16).Value = 1280
$ws.Cells.Item(105, 17).Value = 25
$ws.Cells.Item(105, 18).Value = 'Hortaliza'

# Row 106
$ws.Cells.Item(106, 1).Value = 9
$ws.Cells.Item(106, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(106, 3).Value = 'Metropolitana'
$ws.Cells.Item(106, 4).Value = 44484
$ws.Cells.Item(106, 5).Value = 13
$ws.Cells.Item(106, 6).Value = 100112022
$ws.Cells.Item(106, 7).Value = 'Arveja Verde'
$ws.Cells.Item(106, 8).Value = 'Perfection'
$ws.Cells.Item(106, 9).Value = 'Primera'
$ws.Cells.Item(106, 10).Value = 43
$ws.Cells.Item(106, 11).Value = 22000
$ws.Cells.Item(106, 12).Value = 25000
$ws.Cells.Item(106, 13).Value = 23465
$ws.Cells.Item(106, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(106, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(106, 16).Value = 939
$ws.Cells.Item(106, 17).Value = 25
$ws.Cells.Item(106, 18).Value = 'Hortaliza'

# Row 107
$ws.Cells.Item(107, 1).Value = 9
$ws.Cells.Item(107, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(107, 3).Value = 'Metropolitana'
$ws.Cells.Item(107, 4).Value = 44229
$ws.Cells.Item(107, 5).Value = 13
$ws.Cells.Item(107, 6).Value = 100112022
$ws.Cells.Item(107, 7).Value = 'Arveja Verde'
$ws.Cells.Item(107, 8).Value = 'Sin especificar'
$ws.Cells.Item(107, 9).Value = 'Primera'
$ws.Cells.Item(107, 10).Value = 40
$ws.Cells.Item(107, 11).Value = 24000
$ws.Cells.Item(107, 12).Value = 24000
$ws.Cells.Item(107, 13).Value = 24000
$ws.Cells.Item(107, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(107, 15).Value = 'Carahue'
$ws.Cells.Item(107, 16).Value = 960
$ws.Cells.Item(107, 17).Value = 25
$ws.Cells.Item(107, 18).Value = 'Hortaliza'

# Row 108
$ws.Cells.Item(108, 1).Value = 9
$ws.Cells.Item(108, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(108, 3).Value = 'Metropolitana'
$ws.Cells.Item(108, 4).Value = 44488
$ws.Cells.Item(108, 5).Value = 13
$ws.Cells.Item(108, 6).Value = 100112022
$ws.Cells.Item(108, 7).Value = 'Arveja Verde'
$ws.Cells.Item(108, 8).Value = 'Perfection'
$ws.Cells.Item(108, 9).Value = 'Primera'
$ws.Cells.Item(108, 10).Value = 25
$ws.Cells.Item(108, 11).Value = 24000
$ws.Cells.Item(108, 12).Value = 25000
$ws.Cells.Item(108, 13).Value = 24480
$ws.Cells.Item(108, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(108, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(108, 16).Value = 979
$ws.Cells.Item(108, 17).Value = 25
$ws.Cells.Item(108, 18).Value = 'Hortaliza'

# Row 109
$ws.Cells.Item(109, 1).Value = 9
$ws.Cells.Item(109, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(109, 3).Value = 'Metropolitana'
$ws.Cells.Item(109, 4).Value = 44196
$ws.Cells.Item(109, 5).Value = 13
$ws.Cells.Item(109, 6).Value = 100112022
$ws.Cells.Item(109, 7).Value = 'Arveja Verde'
$ws.Cells.Item(109, 8).Value = 'Sin especificar'
$ws.Cells.Item(109, 9).Value = 'Primera'
$ws.Cells.Item(109, 10).Value = 30
$ws.Cells.Item(109, 11).Value = 26000
$ws.Cells.Item(109, 12).Value = 27000
$ws.Cells.Item(109, 13).Value = 26500
$ws.Cells.Item(109, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(109, 15).Value = 'Carahue'
$ws.Cells.Item(109, 16).Value = 1060
$ws.Cells.Item(109, 17).Value = 25
$ws.Cells.Item(109, 18).Value = 'Hortaliza'

# Row 110
$ws.Cells.Item(110, 1).Value = 9
$ws.Cells.Item(110, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(110, 3).Value = 'Metropolitana'
$ws.Cells.Item(110, 4).Value = 44769
$ws.Cells.Item(110, 5).Value = 13
$ws.Cells.Item(110, 6).Value = 100112022
$ws.Cells.Item(110, 7).Value = 'Arveja Verde'
$ws.Cells.Item(110, 8).Value = 'Perfection'
$ws.Cells.Item(110, 9).Value = 'Primera'
$ws.Cells.Item(110, 10).Value = 25
$ws.Cells.Item(110, 11).Value = 43000
$ws.Cells.Item(110, 12).Value = 43000
$ws.Cells.Item(110, 13).Value = 43000
$ws.Cells.Item(110, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(110, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(110, 16).Value = 1720
$ws.Cells.Item(110, 17).Value = 25
$ws.Cells.Item(110, 18).Value = 'Hortaliza'

# Row 111
$ws.Cells.Item(111, 1).Value = 9
$ws.Cells.Item(111, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(111, 3).Value = 'Metropolitana'
$ws.Cells.Item(111, 4).Value = 44596
$ws.Cells.Item(111, 5).Value = 13
$ws.Cells.Item(111, 6).Value = 100112022
$ws.Cells.Item(111, 7).Value = 'Arveja Verde'
$ws.Cells.Item(111, 8).Value = 'Sin especificar'
$ws.Cells.Item(111, 9).Value = 'Primera'
$ws.Cells.Item(111, 10).Value = 30
$ws.Cells.Item(111, 11).Value = 28000
$ws.Cells.Item(111, 12).Value = 28000
$ws.Cells.Item(111, 13).Value = 28000
$ws.Cells.Item(111, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(111, 15).Value = 'Carahue'
$ws.Cells.Item(111, 16).Value = 1120
$ws.Cells.Item(111, 17).Value = 25
$ws.Cells.Item(111, 18).Value = 'Hortaliza'

# Row 112
$ws.Cells.Item(112, 1).Value = 9
$ws.Cells.Item(112, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(112, 3).Value = 'Metropolitana'
$ws.Cells.Item(112, 4).Value = 44496
$ws.Cells.Item(112, 5).Value = 13
$ws.Cells.Item(112, 6).Value = 100112022
$ws.Cells.Item(112, 7).Value = 'Arveja Verde'
$ws.Cells.Item(112, 8).Value = 'Perfection'
$ws.Cells.Item(112, 9).Value = 'Primera'
$ws.Cells.Item(112, 10).Value = 30
$ws.Cells.Item(112, 11).Value = 24000
$ws.Cells.Item(112, 12).Value = 24000
$ws.Cells.Item(112, 13).Value = 24000
$ws.Cells.Item(112, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(112, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(112, 16).Value = 960
$ws.Cells.Item(112, 17).Value = 25
$ws.Cells.Item(112, 18).Value = 'Hortaliza'

# Row 113
$ws.Cells.Item(113, 1).Value = 9
$ws.Cells.Item(113, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(113, 3).Value = 'Metropolitana'
$ws.Cells.Item(113, 4).Value = 44425
$ws.Cells.Item(113, 5).Value = 13
$ws.Cells.Item(113, 6).Value = 100112022
$ws.Cells.Item(113, 7).Value = 'Arveja Verde'
$ws.Cells.Item(113, 8).Value = 'Perfection'
$ws.Cells.Item(113, 9).Value = 'Primera'
$ws.Cells.Item(113, 10).Value = 16
$ws.Cells.Item(113, 11).Value = 35000
$ws.Cells.Item(113, 12).Value = 36000
$ws.Cells.Item(113, 13).Value = 35500
$ws.Cells.Item(113, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(113, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(113, 16).Value = 1420
$ws.Cells.Item(113, 17).Value = 25
$ws.Cells.Item(113, 18).Value = 'Hortaliza'

# Row 114
$ws.Cells.Item(114, 1).Value = 9
$ws.Cells.Item(114, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(114, 3).Value = 'Metropolitana'
$ws.Cells.Item(114, 4).Value = 44377
$ws.Cells.Item(114, 5).Value = 13
$ws.Cells.Item(114, 6).Value = 100112022
$ws.Cells.Item(114, 7).Value = 'Arveja Verde'
$ws.Cells.Item(114, 8).Value = 'Perfection'
$ws.Cells.Item(114, 9).Value = 'Primera'
$ws.Cells.Item(114, 10).Value = 25
$ws.Cells.Item(114, 11).Value = 39000
$ws.Cells.Item(114, 12).Value = 40000
$ws.Cells.Item(114, 13).Value = 39480
$ws.Cells.Item(114, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(114, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(114, 16).Value = 1579
$ws.Cells.Item(114, 17).Value = 25
$ws.Cells.Item(114, 18).Value = 'Hortaliza'

# Row 115
$ws.Cells.Item(115, 1).Value = 9
$ws.Cells.Item(115, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(115, 3).Value = 'Metropolitana'
$ws.Cells.Item(115, 4).Value = 44512
$ws.Cells.Item(115, 5).Value = 13
$ws.Cells.Item(115, 6).Value = 100112022
$ws.Cells.Item(115, 7).Value = 'Arveja Verde'
$ws.Cells.Item(115, 8).Value = 'Sin especificar'
$ws.Cells.Item(115, 9).Value = 'Primera'
$ws.Cells.Item(115, 10).Value = 34
$ws.Cells.Item(115, 11).Value = 14000
$ws.Cells.Item(115, 12).Value = 15000
$ws.Cells.Item(115, 13).Value = 14500
$ws.Cells.Item(115, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(115, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(115, 16).Value = 580
$ws.Cells.Item(115, 17).Value = 25
$ws.Cells.Item(115, 18).Value = 'Hortaliza'

# Row 116
$ws.Cells.Item(116, 1).Value = 9
$ws.Cells.Item(116, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(116, 3).Value = 'Metropolitana'
$ws.Cells.Item(116, 4).Value = 44181
$ws.Cells.Item(116, 5).Value = 13
$ws.Cells.Item(116, 6).Value = 100112022
$ws.Cells.Item(116, 7).Value = 'Arveja Verde'
$ws.Cells.Item(116, 8).Value = 'Sin especificar'
$ws.Cells.Item(116, 9).Value = 'Primera'
$ws.Cells.Item(116, 10).Value = 25
$ws.Cells.Item(116, 11).Value = 24000
$ws.Cells.Item(116, 12).Value = 24000
$ws.Cells.Item(116, 13).Value = 24000
$ws.Cells.Item(116, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(116, 15).Value = 'Carahue'
$ws.Cells.Item(116, 16).Value = 960
$ws.Cells.Item(116, 17).Value = 25
$ws.Cells.Item(116, 18).Value = 'Hortaliza'

# Row 117
$ws.Cells.Item(117, 1).Value = 9
$ws.Cells.Item(117, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(117, 3).Value = 'Metropolitana'
$ws.Cells.Item(117, 4).Value = 44497
$ws.Cells.Item(117, 5).Value = 13
$ws.Cells.Item(117, 6).Value = 100112022
$ws.Cells.Item(117, 7).Value = 'Arveja Verde'
$ws.Cells.Item(117, 8).Value = 'Sin especificar'
$ws.Cells.Item(117, 9).Value = 'Primera'
$ws.Cells.Item(117, 10).Value = 35
$ws.Cells.Item(117, 11).Value = 1300
$ws.Cells.Item(117, 12).Value = 1500
$ws.Cells.Item(117, 13).Value = 1414
$ws.Cells.Item(117, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(117, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(117, 16).Value = 57
$ws.Cells.Item(117, 17).Value = 25
$ws.Cells.Item(117, 18).Value = 'Hortaliza'

# Row 118
$ws.Cells.Item(118, 1).Value = 9
$ws.Cells.Item(118, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(118, 3).Value = 'Metropolitana'
$ws.Cells.Item(118, 4).Value = 44362
$ws.Cells.Item(118, 5).Value = 13
$ws.Cells.Item(118, 6).Value = 100112022
$ws.Cells.Item(118, 7).Value = 'Arveja Verde'
$ws.Cells.Item(118, 8).Value = 'Perfection'
$ws.Cells.Item(118, 9).Value = 'Primera'
$ws.Cells.Item(118, 10).Value = 16
$ws.Cells.Item(118, 11).Value = 35000
$ws.Cells.Item(118, 12).Value = 37000
$ws.Cells.Item(118, 13).Value = 36000
$ws.Cells.Item(118, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(118, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(118, 16).Value = 1440
$ws.Cells.Item(118, 17).Value = 25
$ws.Cells.Item(118, 18).Value = 'Hortaliza'

# Row 119
$ws.Cells.Item(119, 1).Value = 9
$ws.Cells.Item(119, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(119, 3).Value = 'Metropolitana'
$ws.Cells.Item(119, 4).Value = 44547
$ws.Cells.Item(119, 5).Value = 13
$ws.Cells.Item(119, 6).Value = 100112022
$ws.Cells.Item(119, 7).Value = 'Arveja Verde'
$ws.Cells.Item(119, 8).Value = 'Sin especificar'
$ws.Cells.Item(119, 9).Value = 'Primera'
$ws.Cells.Item(119, 10).Value = 43
$ws.Cells.Item(119, 11).Value = 11000
$ws.Cells.Item(119, 12).Value = 12000
$ws.Cells.Item(119, 13).Value = 11512
$ws.Cells.Item(119, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(119, 15).Value = 'Carahue'
$ws.Cells.Item(119, 16).Value = 460
$ws.Cells.Item(119, 17).Value = 25
$ws.Cells.Item(119, 18).Value = 'Hortaliza'

# Row 120
$ws.Cells.Item(120, 1).Value = 9
$ws.Cells.Item(120, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(120, 3).Value = 'Metropolitana'
$ws.Cells.Item(120, 4).Value = 44657
$ws.Cells.Item(120, 5).Value = 13
$ws.Cells.Item(120, 6).Value = 100112022
$ws.Cells.Item(120, 7).Value = 'Arveja Verde'
$ws.Cells.Item(120, 8).Value = 'Sin especificar'
$ws.Cells.Item(120, 9).Value = 'Primera'
$ws.Cells.Item(120, 10).Value = 25
$ws.Cells.Item(120, 11).Value = 21000
$ws.Cells.Item(120, 12).Value = 21000
$ws.Cells.Item(120, 13).Value = 21000
$ws.Cells.Item(120, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(120, 15).Value = 'Carahue'
$ws.Cells.Item(120, 16).Value = 840
$ws.Cells.Item(120, 17).Value = 25
$ws.Cells.Item(120, 18).Value = 'Hortaliza'

# Row 121
$ws.Cells.Item(121, 1).Value = 9
$ws.Cells.Item(121, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(121, 3).Value = 'Metropolitana'
$ws.Cells.Item(121, 4).Value = 44189
$ws.Cells.Item(121, 5).Value = 13
$ws.Cells.Item(121, 6).Value = 100112022
$ws.Cells.Item(121, 7).Value = 'Arveja Verde'
$ws.Cells.Item(121, 8).Value = 'Sin especificar'
$ws.Cells.Item(121, 9).Value = 'Primera'
$ws.Cells.Item(121, 10).Value = 25
$ws.Cells.Item(121, 11).Value = 28000
$ws.Cells.Item(121, 12).Value = 28000
$ws.Cells.Item(121, 13).Value = 28000
$ws.Cells.Item(121, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(121, 15).Value = 'Carahue'
$ws.Cells.Item(121, 16).Value = 1120
$ws.Cells.Item(121, 17).Value = 25
$ws.Cells.Item(121, 18).Value = 'Hortaliza'
$ws.Cells.Item(121, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 122
$ws.Cells.Item(122, 1).Value = 9
$ws.Cells.Item(122, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(122, 3).Value = 'Metropolitana'
$ws.Cells.Item(122, 4).Value = 44489
$ws.Cells.Item(122, 5).Value = 13
$ws.Cells.Item(122, 6).Value = 100112022
$ws.Cells.Item(122, 7).Value = 'Arveja Verde'
$ws.Cells.Item(122, 8).Value = 'Perfection'
$ws.Cells.Item(122, 9).Value = 'Primera'
$ws.Cells.Item(122, 10).Value = 16
$ws.Cells.Item(122, 11).Value = 24000
$ws.Cells.Item(122, 12).Value = 25000
$ws.Cells.Item(122, 13).Value = 24500
$ws.Cells.Item(122, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(122, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(122, 16).Value = 980
$ws.Cells.Item(122, 17).Value = 25
$ws.Cells.Item(122, 18).Value = 'Hortaliza'
$ws.Cells.Item(122, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
